$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 380
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 380
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = 380
$ws.Range("N12").Value = -720

$ws.Range("H17").Value = 1189.2
$ws.Range("J17").Value = 1189.2
$ws.Range("L17").Value = 3567.6
$ws.Range("N17").Value = -3903.6

$ws.Range("H70").Value = 1540
$ws.Range("I70").Value = 1475
$ws.Range("J70").Value = 1600
$ws.Range("K70").Value = 4425
$ws.Range("L70").Value = 4800
$ws.Range("M70").Value = -4155
$ws.Range("N70").Value = -5340

$ws.Range("H73").Value = 1540
$ws.Range("I73").Value = 1475
$ws.Range("J73").Value = 1600
$ws.Range("K73").Value = 4425
$ws.Range("L73").Value = 4800
$ws.Range("M73").Value = -3489
$ws.Range("N73").Value = -6672

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = ""
$ws.Range("N112").Value = 0

$ws.Range("H113").Value = 4333
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -11507.5

$ws.Range("H135").Value = 1612.8235
$ws.Range("I135").Value = 1447.8667
$ws.Range("K135").Value = 13030.8003
$ws.Range("M135").Value = -10495.8003

$ws.Range("H137").Value = 2185.2856
$ws.Range("I137").Value = 2185.2856
$ws.Range("K137").Value = 6555.8568
$ws.Range("M137").Value = -4005.8568

$ws.Range("H138").Value = 2432.4167
$ws.Range("I138").Value = 733.3333
$ws.Range("J138").Value = 2998.7778
$ws.Range("K138").Value = 2199.9999
$ws.Range("L138").Value = 8996.3334
$ws.Range("M138").Value = 2940.0001
$ws.Range("N138").Value = -19276.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 488
$ws.Range("I2").Value = 488
$ws.Range("K2").Value = 488
$ws.Range("M2").Value = -375

$ws.Range("H26").Value = 196.5
$ws.Range("I26").Value = 196.5
$ws.Range("K26").Value = 196.5
$ws.Range("M26").Value = 133.5

$ws.Range("H92").Value = 74500
$ws.Range("J92").Value = 74500
$ws.Range("L92").Value = 74500
$ws.Range("N92").Value = -79492

$ws.Range("H116").Value = 488
$ws.Range("I116").Value = 488
$ws.Range("K116").Value = 488
$ws.Range("M116").Value = 1806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 488
$ws.Range("I3").Value = 488
$ws.Range("K3").Value = 488
$ws.Range("M3").Value = -374

$ws.Range("H44").Value = 12500
$ws.Range("I44").Value = 12500
$ws.Range("K44").Value = 12500
$ws.Range("M44").Value = -12003

$ws.Range("H92").Value = 24700.5
$ws.Range("J92").Value = 24700.5
$ws.Range("L92").Value = 24700.5
$ws.Range("N92").Value = -29692.5

$ws.Range("H134").Value = 1244
$ws.Range("I134").Value = 1244
$ws.Range("K134").Value = 3732
$ws.Range("M134").Value = -1197

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 100000740
$ws.Range("I12").Value = 100000740
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 100000740
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = -100000570

$ws.Range("H31").Value = 3812
$ws.Range("I31").Value = 2997.2856
$ws.Range("J31").Value = 4952.6
$ws.Range("K31").Value = 2997.2856
$ws.Range("L31").Value = 4952.6
$ws.Range("M31").Value = -2702.2856
$ws.Range("N31").Value = -5542.6

$ws.Range("H34").Value = 3812
$ws.Range("I34").Value = 2997.2856
$ws.Range("J34").Value = 4952.6
$ws.Range("K34").Value = 2997.2856
$ws.Range("L34").Value = 4952.6
$ws.Range("M34").Value = -2795.2856
$ws.Range("N34").Value = -5356.6

$ws.Range("H134").Value = 2798.6
$ws.Range("I134").Value = 2798.6
$ws.Range("K134").Value = 8395.799999999999
$ws.Range("M134").Value = -5860.799999999999

$ws.Range("H141").Value = 559719.5
$ws.Range("J141").Value = 559719.5
$ws.Range("L141").Value = 559719.5
$ws.Range("N141").Value = -570079.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1798
$ws.Range("I107").Value = 1798
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 5394
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -3474

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 1174
$ws.Range("J92").Value = 1174
$ws.Range("L92").Value = 1174
$ws.Range("N92").Value = -4918

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = 0

$ws.Range("H13").Value = 1006
$ws.Range("I13").Value = 1006
$ws.Range("K13").Value = 1006
$ws.Range("M13").Value = -866

$ws.Range("H22").Value = 812.5
$ws.Range("I22").Value = 625
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 625
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -330
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 812.5
$ws.Range("I27").Value = 625
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 625
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -518
$ws.Range("N27").Value = -1214

$ws.Range("H46").Value = 847.5
$ws.Range("I46").Value = 795
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 795
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -607
$ws.Range("N46").Value = -1276

$ws.Range("H55").Value = 541.55554
$ws.Range("I55").Value = 476.66666
$ws.Range("K55").Value = 476.66666
$ws.Range("M55").Value = -303.66666

$ws.Range("H82").Value = 1213.6666
$ws.Range("I82").Value = 961
$ws.Range("J82").Value = 1466.3334
$ws.Range("K82").Value = 961
$ws.Range("L82").Value = 1466.3334
$ws.Range("M82").Value = -600
$ws.Range("N82").Value = -2188.3334

$ws.Range("H85").Value = 1213.6666
$ws.Range("I85").Value = 961
$ws.Range("J85").Value = 1466.3334
$ws.Range("K85").Value = 961
$ws.Range("L85").Value = 1466.3334
$ws.Range("M85").Value = 287
$ws.Range("N85").Value = -3962.3334

$ws.Range("H98").Value = 56460.5
$ws.Range("J98").Value = 56460.5
$ws.Range("L98").Value = 56460.5
$ws.Range("N98").Value = -62450.5

$ws.Range("H101").Value = 18453.666
$ws.Range("J101").Value = 18453.666
$ws.Range("L101").Value = 18453.666
$ws.Range("N101").Value = -24943.666

$ws.Range("H136").Value = 54811.43
$ws.Range("I136").Value = 22592.834
$ws.Range("K136").Value = 67778.50199999999
$ws.Range("M136").Value = -65228.50199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 99999
$ws.Range("J40").Value = 99999
$ws.Range("L40").Value = 99999
$ws.Range("N40").Value = -100297

$ws.Range("H58").Value = 23663
$ws.Range("I58").Value = 23663
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 23663
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = ""
$ws.Range("N58").Value = -23355

$ws.Range("H107").Value = 1393
$ws.Range("I107").Value = 815.8
$ws.Range("K107").Value = 2447.4
$ws.Range("M107").Value = -527.3999999999996

$ws.Range("H113").Value = 497.85715
$ws.Range("I113").Value = 398
$ws.Range("J113").Value = 747.5
$ws.Range("K113").Value = 1194
$ws.Range("L113").Value = 2242.5
$ws.Range("M113").Value = 976
$ws.Range("N113").Value = -6582.5

$ws.Range("H122").Value = 672
$ws.Range("I122").Value = 672
$ws.Range("K122").Value = 2016
$ws.Range("M122").Value = 434

$ws.Range("H124").Value = 92500
$ws.Range("J124").Value = 92500
$ws.Range("L124").Value = 92500
$ws.Range("N124").Value = -102320

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = ""
$ws.Range("N141").Value = 0
